# Update crypto price / volume(1h) figures (and a few re-ranked coin rows)
# on Sheet1 to match the latest scrape.
# Note: some "Price" values look like plain decimal numbers (e.g. "1.00",
# "0.585"); a leading apostrophe is used to force them to stay text cells
# (matching the original inlineStr cells) instead of being auto-converted
# to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.980.58'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.937.68'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''375.43'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').Value = '''101.50'
$ws.Range('E6').Value = '  -3.55%  '
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.585'
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('D10').Value = '''36.24'
$ws.Range('E10').Value = '  -2.81%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = '''0.0849'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').Value = '3.408.07'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = '''17.97'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '''7.41'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').Value = '2.940.94'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '''11.17'
$ws.Range('E17').Value = '  +51.29%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').Value = '''0.979'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = '50.967.80'
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('D20').Value = '''3.14'
$ws.Range('E20').Value = '  -5.75%  '
$ws.Range('D21').Value = '''12.54'
$ws.Range('E21').Value = '  -2.70%  '
$ws.Range('D22').Value = '0.0₃0956'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').Value = '''264.02'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('D24').Value = '''68.53'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').Value = '''3.16'
$ws.Range('E25').Value = '  +11.71%  '
$ws.Range('D26').Value = '''8.06'
$ws.Range('E26').Value = '  +4.15%  '
$ws.Range('D27').Value = '''7.73'
$ws.Range('E27').Value = '  +3.51%  '
$ws.Range('D28').Value = '''0.169'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''25.60'
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.110'
$ws.Range('E31').Value = '  -3.24%  '
$ws.Range('D32').Value = '''9.91'
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('D33').Value = '''50.63'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').Value = '''2.02'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '''32.92'
$ws.Range('E35').Value = '  -6.11%  '
$ws.Range('D36').Value = '''0.0441'
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '''3.11'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''2.51'
$ws.Range('E40').Value = '  -3.53%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '''16.23'
$ws.Range('E41').Value = '  -5.95%  '
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('D43').Value = '''120.78'
$ws.Range('E43').Value = '  -2.94%  '
$ws.Range('D44').Value = '''20.97'
$ws.Range('E44').Value = '  -5.02%  '
$ws.Range('D45').Value = '''0.279'
$ws.Range('E45').Value = '  -4.29%  '
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('D48').Value = '''2.31'
$ws.Range('E48').Value = '  -3.33%  '
$ws.Range('D49').Value = '1.997.59'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('D50').Value = '''0.0335'
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('D51').Value = '''1.29'
$ws.Range('E51').Value = '  +0.15%  '
